# Append newly received loss-of-sale records to the Walk-In Report sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows (row 22 .. row 26), matching the existing table schema:
# #, Date, Customer Name, Contact, Function Date, Staff, Status, Category,
# Sub Category, Repeat count, Remarks
$newRows = @(
    @{ Row = 22; Num = 20; Date = "22-12-2025"; Name = "MAJID";  Contact = 8606416639; FuncDate = "25-12-2025"; Staff = "Vishnu Hari C";     Status = "Loss"; Category = "ENQUIRY"; SubCategory = "Enquiry for Relative/Friend";  Repeat = "-"; Remarks = "LOSS" }
    @{ Row = 23; Num = 21; Date = "23-12-2025"; Name = "fariz";  Contact = 8943665555; FuncDate = "03-01-2026"; Staff = "Siyad vallikkadan";  Status = "Loss"; Category = "ENQUIRY"; SubCategory = "ENQUIRY WITHOUT BRIDE/FAMILY"; Repeat = "-"; Remarks = "tomorrow confirm" }
    @{ Row = 24; Num = 22; Date = "23-12-2025"; Name = "SUHAIL"; Contact = 9447843666; FuncDate = "27-12-2025"; Staff = "Vishnu Hari C";     Status = "Loss"; Category = "PRODUCT"; SubCategory = "PRODUCT NOT AVAILABLE";        Repeat = "-"; Remarks = "LOSS" }
    @{ Row = 25; Num = 23; Date = "24-12-2025"; Name = "MUNEES"; Contact = 9048208991; FuncDate = "27-12-2025"; Staff = "Vishnu Hari C";     Status = "Loss"; Category = "ENQUIRY"; SubCategory = "ENQUIRY WITHOUT BRIDE/FAMILY"; Repeat = "-"; Remarks = "LOSS" }
    @{ Row = 26; Num = 24; Date = "24-12-2025"; Name = "SADHIK"; Contact = 8848446486; FuncDate = "03-01-2026"; Staff = "Siyad vallikkadan";  Status = "Loss"; Category = "ENQUIRY"; SubCategory = "ENQUIRY WITHOUT BRIDE/FAMILY"; Repeat = "-"; Remarks = "LOSS" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Column A (#) - numeric serial, same integer number format used by the
    # rest of the table's "#" column.
    $ws.Cells.Item($row, 1).Value = $r.Num
    $ws.Cells.Item($row, 1).NumberFormat = "0"

    # Column B (Date) - plain text, e.g. "22-12-2025".
    $ws.Cells.Item($row, 2).Value = $r.Date

    # Column C (Customer Name)
    $ws.Cells.Item($row, 3).Value = $r.Name

    # Column D (Contact) - numeric phone number, same integer number format
    # used by the rest of the table's "Contact" column.
    $ws.Cells.Item($row, 4).Value = $r.Contact
    $ws.Cells.Item($row, 4).NumberFormat = "0"

    # Column E (Function Date) - plain text. A handful of these values
    # (day <= 12) are ambiguous and would otherwise auto-convert to a real
    # date serial on assignment, so force literal text via a quote-prefix,
    # exactly like typing '03-01-2026 into the cell in Excel.
    if ($r.FuncDate -match '^\d{2}-\d{2}-\d{4}$') {
        $day = [int]($r.FuncDate.Substring(0,2))
        if ($day -le 12) {
            $ws.Cells.Item($row, 5).Value = "'" + $r.FuncDate
        } else {
            $ws.Cells.Item($row, 5).Value = $r.FuncDate
        }
    } else {
        $ws.Cells.Item($row, 5).Value = $r.FuncDate
    }

    # Column F (Staff)
    $ws.Cells.Item($row, 6).Value = $r.Staff

    # Column G (Status)
    $ws.Cells.Item($row, 7).Value = $r.Status

    # Column H (Category)
    $ws.Cells.Item($row, 8).Value = $r.Category

    # Column I (Sub Category)
    $ws.Cells.Item($row, 9).Value = $r.SubCategory

    # Column J (Repeat count)
    $ws.Cells.Item($row, 10).Value = $r.Repeat

    # Column K (Remarks)
    $ws.Cells.Item($row, 11).Value = $r.Remarks
}
